$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.128864169120789
$ws.Range("B1").Value = 2.098086357116699
$ws.Range("C1").Value = 10.04667663574219
$ws.Range("D1").Value = 2.500462055206299
$ws.Range("E1").Value = 1.298213124275208
